$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("G7").Value = 2.1
$ws.Range("I7").Value = 3.2
$ws.Range("L7").Value = 3.6
$ws.Range("Q7").Value = 1.62
$ws.Range("R7").Value = 2.25
$ws.Range("G8").Value = 2.4
$ws.Range("I8").Value = 2.75
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 1.73
$ws.Range("R8").Value = 2.08
$ws.Range("AA8").Value = 19
$ws.Range("AI8").Value = 10
$ws.Range("AZ8").Value = 41
$ws.Range("I10").Value = 2.05
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 2.88
$ws.Range("Y10").Value = 12
$ws.Range("AD10").Value = 6.5
$ws.Range("AI10").Value = 9.5
$ws.Range("AJ10").Value = 19
$ws.Range("AQ10").Value = 67
$ws.Range("AS10").Value = 251
$ws.Range("I14").Value = 4.2
$ws.Range("U14").Value = 1.57
$ws.Range("V14").Value = 2.25
$ws.Range("AF14").Value = 41
$ws.Range("AO14").Value = 9
$ws.Range("N20").Value = 5.9
